$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E8").Value = "1m77"
$ws.Range("E9").Value = "1m78"

$ws.Range("F18").Select()
$excel.ActiveWindow.ScrollRow = 2
